$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab label (the visible sheet name in workbook.xml)
$ws.Name = "Through 2021-12-16"

# Update the December header label to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-16)"

# Update December (row 13) counts
$ws.Range("B13").Value = 19
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = 57
$ws.Range("E13").Value = 35
$ws.Range("F13").Value = 27
$ws.Range("G13").Value = 78
$ws.Range("H13").Value = 125

# Update Total (row 14) counts
$ws.Range("B14").Value = 310
$ws.Range("C14").Value = 613
$ws.Range("D14").Value = 878
$ws.Range("E14").Value = 717
$ws.Range("F14").Value = 561
$ws.Range("G14").Value = 1342
$ws.Range("H14").Value = 1768
